$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 1735.3422
$ws.Range("I15").Value = 1735.3422
$ws.Range("K15").Value = 5206.0266
$ws.Range("M15").Value = -5037.0266
# Row 64
$ws.Range("H64").Value = 6403.222
$ws.Range("I64").Value = 4861.4546
$ws.Range("J64").Value = 8826
$ws.Range("K64").Value = 4861.4546
$ws.Range("L64").Value = 8826
$ws.Range("M64").Value = -4613.4546
$ws.Range("N64").Value = -9322
# Row 67
$ws.Range("H67").Value = 6403.222
$ws.Range("I67").Value = 4861.4546
$ws.Range("J67").Value = 8826
$ws.Range("K67").Value = 4861.4546
$ws.Range("L67").Value = 8826
$ws.Range("M67").Value = -4003.4546
$ws.Range("N67").Value = -10542
# Row 74
$ws.Range("H74").Value = 8983.909
$ws.Range("I74").Value = 7938.5
$ws.Range("K74").Value = 7938.5
$ws.Range("M74").Value = -7002.5
# Row 77
$ws.Range("H77").Value = 8983.909
$ws.Range("I77").Value = 7938.5
$ws.Range("K77").Value = 39692.5
$ws.Range("M77").Value = -35012.5
# Row 100
$ws.Range("H100").Value = 2673.0833
$ws.Range("I100").Value = 2083.5454
$ws.Range("K100").Value = 2083.5454
$ws.Range("M100").Value = -1542.5454
# Row 112
$ws.Range("H112").Value = 2171.5
$ws.Range("J112").Value = 2205.8
$ws.Range("L112").Value = 6617.400000000001
$ws.Range("N112").Value = -8833.400000000001
# Row 121
$ws.Range("H121").Value = 2271.8572
$ws.Range("J121").Value = 2317.6667
$ws.Range("L121").Value = 6953.000100000001
$ws.Range("N121").Value = -10447.0001
# Row 131
$ws.Range("H131").Value = 2853.0908
$ws.Range("I131").Value = 2738.4
$ws.Range("K131").Value = 8215.200000000001
$ws.Range("M131").Value = -3175.200000000001
# Row 137
$ws.Range("H137").Value = 31255758
$ws.Range("I137").Value = 45461364
$ws.Range("K137").Value = 136384092
$ws.Range("M137").Value = -136381542
# Row 138
$ws.Range("H138").Value = 8374.329
$ws.Range("I138").Value = 5278.3335
$ws.Range("J138").Value = 8600.866
$ws.Range("K138").Value = 15835.0005
$ws.Range("L138").Value = 25802.598
$ws.Range("M138").Value = -10695.0005
$ws.Range("N138").Value = -36082.598
# Row 141
$ws.Range("H141").Value = 7391
$ws.Range("I141").Value = 4556.136
$ws.Range("J141").Value = 11845.786
$ws.Range("K141").Value = 13668.408
$ws.Range("L141").Value = 35537.358
$ws.Range("M141").Value = -8488.408000000001
$ws.Range("N141").Value = -45897.358

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 16077.113
$ws.Range("I32").Value = 7819.4517
$ws.Range("J32").Value = 35768.46
$ws.Range("K32").Value = 7819.4517
$ws.Range("L32").Value = 35768.46
$ws.Range("M32").Value = -7532.4517
$ws.Range("N32").Value = -36342.46
# Row 45
$ws.Range("H45").Value = 1291.75
$ws.Range("I45").Value = 1136.4
$ws.Range("K45").Value = 1136.4
$ws.Range("M45").Value = -759.4000000000001
# Row 61
$ws.Range("H61").Value = 4452
$ws.Range("I61").Value = 4025.577
$ws.Range("J61").Value = 9995.5
$ws.Range("K61").Value = 4025.577
$ws.Range("L61").Value = 9995.5
$ws.Range("M61").Value = -3813.577
$ws.Range("N61").Value = -10419.5
# Row 122
$ws.Range("H122").Value = 6674.2
$ws.Range("J122").Value = 6178.25
$ws.Range("L122").Value = 18534.75
$ws.Range("N122").Value = -23434.75
# Row 132
$ws.Range("H132").Value = 15613.19
$ws.Range("I132").Value = 11928.406
$ws.Range("J132").Value = 27404.5
$ws.Range("K132").Value = 35785.218
$ws.Range("L132").Value = 82213.5
$ws.Range("M132").Value = -33255.218
$ws.Range("N132").Value = -87273.5
# Row 136
$ws.Range("H136").Value = 4452
$ws.Range("I136").Value = 4025.577
$ws.Range("J136").Value = 9995.5
$ws.Range("K136").Value = 12076.731
$ws.Range("L136").Value = 29986.5
$ws.Range("M136").Value = -9526.731
$ws.Range("N136").Value = -35086.5

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 4553.5557
$ws.Range("I20").Value = 4033
$ws.Range("K20").Value = 4033
$ws.Range("M20").Value = -3786
# Row 88
$ws.Range("H88").Value = 23153.143
$ws.Range("J88").Value = 22325
$ws.Range("L88").Value = 22325
$ws.Range("N88").Value = -23137
# Row 91
$ws.Range("H91").Value = 23153.143
$ws.Range("J91").Value = 22325
$ws.Range("L91").Value = 22325
$ws.Range("N91").Value = -25133
# Row 134
$ws.Range("H134").Value = 3739.6
$ws.Range("I134").Value = 1566
$ws.Range("K134").Value = 4698
$ws.Range("M134").Value = -2163
# Row 140
$ws.Range("H140").Value = 86500
$ws.Range("J140").Value = 94653.84
$ws.Range("L140").Value = 94653.84
$ws.Range("N140").Value = -105013.84

$ws = $wb.Worksheets.Item("CRP")
# Row 94
$ws.Range("H94").Value = 1514.6428
$ws.Range("I94").Value = 1609.8334
$ws.Range("J94").Value = 1443.25
$ws.Range("K94").Value = 1609.8334
$ws.Range("L94").Value = 1443.25
$ws.Range("M94").Value = -1158.8334
$ws.Range("N94").Value = -2345.25
# Row 99
$ws.Range("H99").Value = 25400710
$ws.Range("I99").Value = 10530451
$ws.Range("K99").Value = 10530451
$ws.Range("M99").Value = -10528953
# Row 122
$ws.Range("H122").Value = 2925
$ws.Range("I122").Value = 2908.875
$ws.Range("K122").Value = 8726.625
$ws.Range("M122").Value = -6276.625
# Row 126
$ws.Range("H126").Value = 25400710
$ws.Range("I126").Value = 10530451
$ws.Range("K126").Value = 31591353
$ws.Range("M126").Value = -31588883
# Row 132
$ws.Range("H132").Value = 2245.6216
$ws.Range("I132").Value = 2056.5
$ws.Range("J132").Value = 3056.1428
$ws.Range("K132").Value = 6169.5
$ws.Range("L132").Value = 9168.428400000001
$ws.Range("M132").Value = -3639.5
$ws.Range("N132").Value = -14228.4284
# Row 137
$ws.Range("H137").Value = 69000
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 1422.5
$ws.Range("I34").Value = 1422.5
$ws.Range("K34").Value = 4267.5
$ws.Range("M34").Value = -4183.5
# Row 39
$ws.Range("H39").Value = 7249.75
$ws.Range("J39").Value = 9000
$ws.Range("L39").Value = 27000
$ws.Range("N39").Value = -27588
# Row 137
$ws.Range("H137").Value = 6921.875
$ws.Range("J137").Value = 9250
$ws.Range("L137").Value = 27750
$ws.Range("N137").Value = -37950

$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 13050
$ws.Range("I5").Value = 9000
$ws.Range("J5").Value = 14400
$ws.Range("K5").Value = 9000
$ws.Range("L5").Value = 14400
$ws.Range("M5").Value = -8888
$ws.Range("N5").Value = -14624
# Row 24
$ws.Range("H24").Value = 16799.8
$ws.Range("J24").Value = 16799.8
$ws.Range("L24").Value = 16799.8
$ws.Range("N24").Value = -17145.8
# Row 80
$ws.Range("H80").Value = 5499.75
$ws.Range("J80").Value = 7666.3335
$ws.Range("L80").Value = 7666.3335
$ws.Range("N80").Value = -9662.333500000001
# Row 83
$ws.Range("H83").Value = 5499.75
$ws.Range("J83").Value = 7666.3335
$ws.Range("L83").Value = 38331.6675
$ws.Range("N83").Value = -48315.6675
# Row 102
$ws.Range("H102").Value = 2781.3784
$ws.Range("I102").Value = 2035.6207
$ws.Range("K102").Value = 2035.6207
$ws.Range("M102").Value = -413.6206999999999
# Row 122
$ws.Range("H122").Value = 4295.9375
$ws.Range("I122").Value = 4056.6155
$ws.Range("K122").Value = 12169.8465
$ws.Range("M122").Value = -9719.8465
# Row 126
$ws.Range("H126").Value = 5062.485
$ws.Range("I126").Value = 5021.269
$ws.Range("J126").Value = 5215.5713
$ws.Range("K126").Value = 15063.807
$ws.Range("L126").Value = 15646.7139
$ws.Range("M126").Value = -12593.807
$ws.Range("N126").Value = -20586.7139
# Row 132
$ws.Range("H132").Value = 5145.963
$ws.Range("I132").Value = 5166.4546
$ws.Range("K132").Value = 15499.3638
$ws.Range("M132").Value = -12969.3638

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 71435860
$ws.Range("I7").Value = 83340170
$ws.Range("K7").Value = 83340170
$ws.Range("M7").Value = -83340058
# Row 40
$ws.Range("H40").Value = 41672400
$ws.Range("I40").Value = 41672400
$ws.Range("K40").Value = 41672400
$ws.Range("M40").Value = -41672264
# Row 61
$ws.Range("H61").Value = 4411.4546
$ws.Range("I61").Value = 4554.857
$ws.Range("K61").Value = 4554.857
$ws.Range("M61").Value = -4352.857
# Row 113
$ws.Range("H113").Value = 4411.4546
$ws.Range("I113").Value = 4554.857
$ws.Range("K113").Value = 4554.857
$ws.Range("M113").Value = -2384.857
# Row 122
$ws.Range("H122").Value = 9054.888999999999
$ws.Range("I122").Value = 6168.3076
$ws.Range("K122").Value = 18504.9228
$ws.Range("M122").Value = -16054.9228
# Row 126
$ws.Range("H126").Value = 71435860
$ws.Range("I126").Value = 83340170
$ws.Range("K126").Value = 250020510
$ws.Range("M126").Value = -250018040
# Row 132
$ws.Range("H132").Value = 8331.200000000001
$ws.Range("I132").Value = 8145.9443
$ws.Range("K132").Value = 24437.8329
$ws.Range("M132").Value = -21907.8329
# Row 136
$ws.Range("H136").Value = 7441.0527
$ws.Range("I136").Value = 7459.067
$ws.Range("K136").Value = 22377.201
$ws.Range("M136").Value = -19827.201

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 53999
$ws.Range("I2").Value = 28000
$ws.Range("K2").Value = 28000
$ws.Range("M2").Value = -27888
# Row 122
$ws.Range("H122").Value = 10468355
$ws.Range("I122").Value = 6253064.5
$ws.Range("K122").Value = 18759193.5
$ws.Range("M122").Value = -18756743.5

